$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.447.48"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "'1.841.00"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'261.18"
$ws.Range("E5").Value = "  -5.86%  "

$ws.Range("D7").Value = "'0.5189"
$ws.Range("E7").Value = "  -1.89%  "

$ws.Range("D8").Value = "'0.3278"
$ws.Range("E8").Value = "  -4.17%  "

$ws.Range("D9").Value = "'0.06789"
$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("D10").Value = "'18.66"
$ws.Range("E10").Value = "  -6.74%  "

$ws.Range("D11").Value = "'0.7695"
$ws.Range("E11").Value = "  -4.20%  "

$ws.Range("D12").Value = "'0.07707"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "'1.833.57"
$ws.Range("E13").Value = "  -2.39%  "

$ws.Range("D14").Value = "'88.33"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("D15").Value = "'5.035"
$ws.Range("E15").Value = "  -2.67%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "'13.92"
$ws.Range("E17").Value = "  -4.28%  "

$ws.Range("D18").Value = "'0.000007985"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "'0.9998"

$ws.Range("D20").Value = "'26.454.33"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("D21").Value = "'2.068.54"
$ws.Range("E21").Value = "  -2.60%  "

$ws.Range("D22").Value = "'4.576"

$ws.Range("D23").Value = "'9.472"
$ws.Range("E23").Value = "  -5.48%  "

$ws.Range("D24").Value = "'5.971"
$ws.Range("E24").Value = "  -3.68%  "

$ws.Range("D25").Value = "'144.20"
$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("D26").Value = "'2.215"
$ws.Range("E26").Value = "  -7.57%  "

$ws.Range("D27").Value = "'1.648"
$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("D28").Value = "'17.00"
$ws.Range("E28").Value = "  -1.85%  "

$ws.Range("D29").Value = "'111.55"
$ws.Range("E29").Value = "  -1.72%  "

$ws.Range("D30").Value = "'4.194"
$ws.Range("E30").Value = "  -3.34%  "

$ws.Range("D31").Value = "'4.137"
$ws.Range("E31").Value = "  -3.77%  "

$ws.Range("D32").Value = "'0.08727"
$ws.Range("E32").Value = "  -1.70%  "

$ws.Range("D33").Value = "'0.04805"
$ws.Range("E33").Value = "  -1.69%  "

$ws.Range("D34").Value = "'1.132"
$ws.Range("E34").Value = "  -3.61%  "

$ws.Range("D35").Value = "'2.836"
$ws.Range("E35").Value = "  -2.05%  "

$ws.Range("D36").Value = "'0.7071"
$ws.Range("E36").Value = "  -2.60%  "

$ws.Range("D37").Value = "'3.072"
$ws.Range("E37").Value = "  -6.46%  "

$ws.Range("D38").Value = "'2.224"
$ws.Range("E38").Value = "  -5.13%  "

$ws.Range("D39").Value = "'0.01761"
$ws.Range("E39").Value = "  -4.31%  "

$ws.Range("D40").Value = "'0.4834"
$ws.Range("E40").Value = "  -5.62%  "

$ws.Range("D41").Value = "'111.45"

$ws.Range("D42").Value = "'0.8913"
$ws.Range("E42").Value = "  -6.71%  "

$ws.Range("D43").Value = "'6.078"
$ws.Range("E43").Value = "  -1.54%  "

$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "'7.715"
$ws.Range("E45").Value = "  -4.50%  "

$ws.Range("D46").Value = "'0.4145"

$ws.Range("D47").Value = "'0.05867"
$ws.Range("E47").Value = "  -1.33%  "

$ws.Range("D48").Value = "'9.014"
$ws.Range("E48").Value = "  -3.25%  "

$ws.Range("D49").Value = "'35.02"
$ws.Range("E49").Value = "  -3.14%  "

$ws.Range("D50").Value = "'0.1220"
$ws.Range("E50").Value = "  -8.93%  "

$ws.Range("D51").Value = "'0.8877"
$ws.Range("E51").Value = "  +0.59%  "
